# Word COM-interop script: collapse the split "<id>...</id>" runs into a
# single run and prefix the inner id value with "p" (113r_1 -> p113r_1,
# 113v_1 -> p113v_1), matching the newly downloaded tc/tcn/tl ids.

$d = $word.ActiveDocument

$d.Content.Find.Execute(
    "<id>113r_1</id>", $true, $false, $false, $false, $false,
    $true, 1, $false, "<id>p113r_1</id>", 2
) | Out-Null

$d.Content.Find.Execute(
    "<id>113v_1</id>", $true, $false, $false, $false, $false,
    $true, 1, $false, "<id>p113v_1</id>", 2
) | Out-Null
